# Apply menu updates: fix a few item names and append newly added menu rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Rename a few existing menu items -------------------------------
# "Frango com Brocolis" -> "Frango com Molho de Brocolis" (rows 10-11)
# "Fatia Presunto"      -> "Presunto de Peru"              (row 52)
# "Melão"               -> "Melao"                          (row 59)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($name -eq "Frango com Brocolis") {
        $ws.Cells.Item($r, 1).Value = "Frango com Molho de Brocolis"
    } elseif ($name -eq "Fatia Presunto") {
        $ws.Cells.Item($r, 1).Value = "Presunto de Peru"
    } elseif ($name -eq "Melão") {
        $ws.Cells.Item($r, 1).Value = "Melao"
    }
}

# --- 2. Append the new menu rows ---------------------------------------

$newRows = New-Object 'object[,]' 29,4

$newRows[0,0]="Atum Continente em Azeite";          $newRows[0,1]="Almoço";       $newRows[0,2]="Proteina";    $newRows[0,3]=225
$newRows[1,0]="Atum Continente em Azeite";          $newRows[1,1]="Jantar";       $newRows[1,2]="Proteina";    $newRows[1,3]=225
$newRows[2,0]="Atum Continente em Oleo";            $newRows[2,1]="Almoço";       $newRows[2,2]="Proteina";    $newRows[2,3]=210
$newRows[3,0]="Atum Continente em Oleo";            $newRows[3,1]="Jantar";       $newRows[3,2]="Proteina";    $newRows[3,3]=210
$newRows[4,0]="Alface";                              $newRows[4,1]="Almoço";       $newRows[4,2]="Verduras";    $newRows[4,3]=15
$newRows[5,0]="Alface";                              $newRows[5,1]="Jantar";       $newRows[5,2]="Verduras";    $newRows[5,3]=15
$newRows[6,0]="Iogurte Grego com Limão Auchan";     $newRows[6,1]="Lanche Manhã"; $newRows[6,2]=$null;         $newRows[6,3]=135
$newRows[7,0]="Iogurte Grego com Limão Auchan";     $newRows[7,1]="Lanche Tarde"; $newRows[7,2]=$null;         $newRows[7,3]=135
$newRows[8,0]="Iogurte Grego com Limão Auchan";     $newRows[8,1]="Café da Manhã";$newRows[8,2]=$null;         $newRows[8,3]=135
$newRows[9,0]="Salada de Alface e Atum com molho Cesar"; $newRows[9,1]="Almoço";  $newRows[9,2]="Proteina";    $newRows[9,3]=160
$newRows[10,0]="Salada de Alface e Atum com molho Cesar";$newRows[10,1]="Jantar"; $newRows[10,2]="Proteina";   $newRows[10,3]=160
$newRows[11,0]="Batata Palha Continente";           $newRows[11,1]="Almoço";      $newRows[11,2]="Carboidrato";$newRows[11,3]=486
$newRows[12,0]="Batata Palha Continente";           $newRows[12,1]="Jantar";      $newRows[12,2]="Carboidrato";$newRows[12,3]=486
$newRows[13,0]="Morango";                            $newRows[13,1]="Café da Manhã";$newRows[13,2]=$null;      $newRows[13,3]=33
$newRows[14,0]="Morango";                            $newRows[14,1]="Lanche Manhã"; $newRows[14,2]=$null;      $newRows[14,3]=33
$newRows[15,0]="Morango";                            $newRows[15,1]="Lanche Tarde"; $newRows[15,2]=$null;      $newRows[15,3]=33
$newRows[16,0]="Geleia de Morango Auchan (-30% de Acucares)"; $newRows[16,1]="Café da Manhã"; $newRows[16,2]=$null; $newRows[16,3]=166
$newRows[17,0]="Geleia de Morango Auchan (-30% de Acucares)"; $newRows[17,1]="Lanche Manhã";  $newRows[17,2]=$null; $newRows[17,3]=166
$newRows[18,0]="Geleia de Morango Auchan (-30% de Acucares)"; $newRows[18,1]="Lanche Tarde";  $newRows[18,2]=$null; $newRows[18,3]=166
$newRows[19,0]="Banana";                             $newRows[19,1]="Café da Manhã"; $newRows[19,2]=$null;     $newRows[19,3]=89
$newRows[20,0]="Banana";                             $newRows[20,1]="Lanche Manhã";  $newRows[20,2]=$null;     $newRows[20,3]=89
$newRows[21,0]="Banana";                             $newRows[21,1]="Lanche Tarde";  $newRows[21,2]=$null;     $newRows[21,3]=89
$newRows[22,0]="Chips de Banana Auchan";             $newRows[22,1]="Café da Manhã"; $newRows[22,2]=$null;     $newRows[22,3]=538
$newRows[23,0]="Chips de Banana Auchan";             $newRows[23,1]="Lanche Manhã";  $newRows[23,2]=$null;     $newRows[23,3]=538
$newRows[24,0]="Chips de Banana Auchan";             $newRows[24,1]="Lanche Tarde";  $newRows[24,2]=$null;     $newRows[24,3]=538
$newRows[25,0]="Pizza de Mussarella";                $newRows[25,1]="Jantar";        $newRows[25,2]="Proteina";$newRows[25,3]=330
$newRows[26,0]="Pizza de Presunto com Milho";        $newRows[26,1]="Jantar";        $newRows[26,2]="Proteina";$newRows[26,3]=330
$newRows[27,0]="Filetes de Pescada";                 $newRows[27,1]="Almoço";        $newRows[27,2]="Proteina";$newRows[27,3]=141
$newRows[28,0]="Filetes de Pescada";                 $newRows[28,1]="Jantar";        $newRows[28,2]="Proteina";$newRows[28,3]=141

$startRow = $lastRow + 1
$endRow = $startRow + 28
$target = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 4))
$target.Value = $newRows

# --- 3. Leave the selection on the last appended cell, like the author --
$ws.Range("D" + $endRow).Select() | Out-Null
